# ---------------------------------------------------------------------------
# 0005-mid-ocean_ridge.pptx edit script
#
# Implements (per the target diff):
#   1. Re-cache the "datetimeFigureOut" footer field text on the slide
#      master and every slide layout that has one (7/23/20 -> 1/8/21).
#   2. Nudge the "Mid-ocean Ridge" title textbox on slide 1 down.
#   3. Nudge the braille title textbox on slide 2 down.
#   4. On slide 2, swap three braille caption boxes ("lithosphere / mantle",
#      "lithosphere / mantle" (right), "asthenosphere / mantle" (centered))
#      for three new ones with updated copy ("lithospheric / mantle",
#      "lith. / mantle", "asthenospheric / mantle" (centered)).
# ---------------------------------------------------------------------------

# EMU -> points helper. `iron_native`'s COM shim rounds Width/Height/Left/Top
# down when the value lands exactly on an integer EMU boundary, so nudge by a
# hair (well under half an EMU-in-points) to land on the correct EMU.
function Pt([double]$emu) {
    return ($emu / 12700.0) + 0.00002
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date field re-cache: "7/23/20" -> "1/8/21" on the slide master and every
#    slide layout exposing a Date Placeholder.
# ---------------------------------------------------------------------------
$newDate = "1/8/21"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 1: move the "Mid-ocean Ridge" title down (y: 160892 -> 362911).
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 1") {
        $sh.Top = Pt(362911)
    }
}

# ---------------------------------------------------------------------------
# 3. Slide 2: move the braille title down (y: 128995 -> 245953).
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $sh = $slide2.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 63") {
        $sh.Top = Pt(245953)
    }
}

# ---------------------------------------------------------------------------
# 4. Slide 2: replace the three braille caption boxes.
#    Template shapes ("Rectangle 84" left-aligned body, "Rectangle 88"
#    centered body) are duplicated so the new shapes inherit identical
#    formatting (font, size, <a:lstStyle/>, paragraph alignment), then are
#    repositioned/resized/retexted/renamed; the three originals are deleted
#    afterwards.
# ---------------------------------------------------------------------------
$rect84 = $null
$rect86 = $null
$rect88 = $null
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $sh = $slide2.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 84") { $rect84 = $sh }
    elseif ($sh.Name -eq "Rectangle 86") { $rect86 = $sh }
    elseif ($sh.Name -eq "Rectangle 88") { $rect88 = $sh }
}

# -- new shape 1: "lithospheric mantle" (left aligned, based on Rectangle 84)
$new1 = $rect84.Duplicate().Item(1)
$new1.Name = "Rectangle 1"
$new1.TextFrame.WordWrap = 0
$new1.Left = Pt(262495)
$new1.Top = Pt(4775422)
$new1.Width = Pt(3081293)
$new1.Height = Pt(830997)
$tr1 = $new1.TextFrame.TextRange
$tr1.Paragraphs(1).Text = "*"
$tr1.Paragraphs(1).Text = "⠇⠊⠞⠓⠕⠎⠏⠓⠑⠗⠊⠉⠀"

# -- new shape 2: "lith. mantle" (left aligned, based on Rectangle 84)
$new2 = $rect84.Duplicate().Item(1)
$new2.Name = "Rectangle 2"
$new2.TextFrame.WordWrap = 0
$new2.Left = Pt(5355370)
$new2.Top = Pt(5136930)
$new2.Width = Pt(1579278)
$new2.Height = Pt(830997)
$tr2 = $new2.TextFrame.TextRange
$tr2.Paragraphs(1).Text = "*"
$tr2.Paragraphs(1).Text = "⠇⠊⠞⠓⠲⠀"

# -- new shape 3: "asthenospheric mantle" (centered, based on Rectangle 88)
$new3 = $rect88.Duplicate().Item(1)
$new3.Name = "Rectangle 9"
$new3.TextFrame.WordWrap = 0
$new3.Left = Pt(2035128)
$new3.Top = Pt(5923738)
$new3.Width = Pt(3526928)
$new3.Height = Pt(830997)
$tr3 = $new3.TextFrame.TextRange
$tr3.Paragraphs(1).Text = "*"
$tr3.Paragraphs(1).Text = "⠁⠎⠞⠓⠑⠝⠕⠎⠏⠓⠑⠗⠊⠉⠀"

# -- drop the old caption boxes
$rect84.Delete()
$rect86.Delete()
$rect88.Delete()
